$wb = $excel.ActiveWorkbook

# 1. Reorder worksheets so that "review_info" comes before "hotel_info"
$hotelInfo = $wb.Worksheets.Item("hotel_info")
$reviewInfo = $wb.Worksheets.Item("review_info")
$hotelInfo.Move($null, $reviewInfo)

# 2. Insert a new "State" column into hotel_info, between Hotel_Name and City
$ws = $wb.Worksheets.Item("hotel_info")
$ws.Columns.Item(3).Insert()
$ws.Range("C1").Value = "State"
$ws.Range("C2").Value = "Louisiana"
